# 9.c.1 — add a 2020 data column (L) to the mobile-network-coverage table,
# mirroring the formatting of the neighbouring 2019 column (K), and move the
# selection the author ended the edit on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (thin divider row above the header): extend the bottom border into L3.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Row 4 (year headers): L4 gets the new "2020" header, formatted like K4 (2019).
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2020

# Row 5 (2G coverage): new 2020 value, formatted like the "0.0"-style cells (H5/I5).
$ws.Range("H5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 99

# Row 6 (3G coverage): new 2020 value, formatted like K6.
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 89.3

# Row 7 (4G coverage): new 2020 value, formatted like K7.
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 81.900000000000006

# Leave the view the way the author left it: scrolled right a touch, with N13 active.
$ws.Range("N13").Select()
